# Nalco aluminium ingot price sheet - append today's price row.
# The source table is newest-first: a new row is inserted at row 2 (right
# below the header) carrying the new date "19-11-2025" with the same
# price/description/circular info as the (old) most-recent row, and every
# other row shifts down by one. The previously-last row (old row 105) keeps
# its data but, because there is no older data point to compare it
# against, also ends up duplicated into the new final row 106.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room: insert a blank row right after the header. This shifts
#    rows 2..105 down to 3..106, preserving cell types/styles/values as-is
#    (no re-parsing of text, so date-like strings stay text).
$ws.Rows.Item(2).Insert()

# 2) Populate the new row 2 by duplicating the row right below it (which is
#    the former row 2), then overwrite just the Date cell.
$ws.Range("A3:F3").Copy($ws.Range("A2:F2"))
$ws.Range("A2").Replace("18-11-2025", "19-11-2025")

# 3) Rebuild the hyperlinks collection from scratch. Row-insert/copy in
#    this host does not renumber the <hyperlinks> anchors, so the cleanest
#    fix is to drop them all and re-add one per F-column cell, pointing at
#    whatever URL text is already sitting in that cell - then restore the
#    plain text style PasteSpecial clobbers when it re-applies the
#    built-in Hyperlink style.
$ws.Hyperlinks.Delete()

for ($r = 2; $r -le 106; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $addr = $cell.Value()
    $ws.Hyperlinks.Add($cell, $addr)
    $ws.Range("F3").Copy()
    $cell.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
